$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The averaging code was re-run to add results for the new spiral sampling
# schemes. Three more rows (17-19) are appended below the existing table,
# continuing the existing A-column index / B-column scheme-name pattern,
# with all of the per-reflection averaged-intensity columns (C:M) equal to 1.

$newRows = @(
    @{ Row = 17; Index = 15; Label = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; Index = 16; Label = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; Index = 17; Label = "HexGrid-60degTilt5degRes" }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    $ws.Cells.Item($r, 1).Value = $nr.Index
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)  # xlPasteFormats, matches style of the row above

    $ws.Cells.Item($r, 2).Value = $nr.Label

    for ($c = 3; $c -le 13; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
